{"js": "// Find the paragraph that ends the Mage section (\"Will periodically bombard...\")\n// and insert the new \"Boss:\" section right after it, before the trailing empty paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"Will periodically bombard the field\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the Mage bombard paragraph\");\n}\n\n// Insert the \"Boss:\" heading paragraph right after the found paragraph.\nconst bossHeading = target.insertParagraph(\"Boss:\", Word.InsertLocation.after);\n\n// Insert the boss description paragraph right after the heading.\nbossHeading.insertParagraph(\n  \"You control the movement of the boss with the arrow keys, the boss will automatically attack. The attack pattern is shown at the right side of the screen\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Locate the paragraph that ends the Mage section (\"Will periodically bombard...\")\n# and insert the new \"Boss:\" section right after it, before the trailing empty paragraph.\n$d = $word.ActiveDocument\n\n$searchText = \"Will periodically bombard the field; showing the locations 1 turn before, all heroes drank a fire resistance potion before the fight.\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$searchText*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find the Mage bombard paragraph\"\n}\n\n# Insert the \"Boss:\" heading paragraph right after the found paragraph.\n$target.Range.InsertParagraphAfter()\n$bossHeading = $d.Paragraphs($target.Index + 1)\n$bossHeading.Range.Text = \"Boss:\"\n\n# Insert the boss description paragraph right after the heading.\n$bossHeading.Range.InsertParagraphAfter()\n$bossDescription = $d.Paragraphs($bossHeading.Index + 1)\n$bossDescription.Range.Text = \"You control the movement of the boss with the arrow keys, the boss will automatically attack. The attack pattern is shown at the right side of the screen\"\n"}
